# Auto-update draw results: append the 2025-12-15 Pick 4 draw as row 90.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 90

# Lead with an apostrophe so values that look numeric/date-like ("2025-12-15",
# "251215") are entered as literal text, matching every other row in this
# results table (all cells are plain text, not real dates/numbers).
$ws.Cells.Item($row, 1).Value = "'2025-12-15"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "'251215"
$ws.Cells.Item($row, 4).Value = "9-4-0-1"
$ws.Cells.Item($row, 5).Value = "2025-12-15T21:46:53.491+04:00"

# Drop the implicit "quote prefix" / number formatting the text-entry above
# picks up, so the new row carries no cell style — consistent with the rest
# of the sheet, which has no per-cell styles at all.
$ws.Range("A$row`:E$row").ClearFormats()
